$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Character / skill code name fixes (CharData, SkillData 코드명 수정) ---
# Character code column (B)
$ws.Range("B2").Value = "Char_K_King001"
$ws.Range("B3").Value = "Char_S_Knight001"
$ws.Range("B4").Value = "Char_S_Archer001"
$ws.Range("B5").Value = "Char_S_Magic001"
$ws.Range("B6").Value = "Char_H_Knight001"
$ws.Range("B7").Value = "Char_H_Archer001"
$ws.Range("B8").Value = "Char_H_Magic001"

# Skill code column (T)
$ws.Range("T2").Value = "Skill_K_King001"
$ws.Range("T3").Value = "Skill_S_Knight001"
$ws.Range("T4").Value = "Skill_S_Archer001"
$ws.Range("T5").Value = "Skill_S_Magic001"
$ws.Range("T6").Value = "Skill_H_Knight001"
$ws.Range("T7").Value = "Skill_H_Archer001"
$ws.Range("T8").Value = "Skill_H_Magic001"

# --- Growth-value (성장치) numeric tweaks ---
$ws.Range("H2").Value = 1.5
$ws.Range("L2").Value = 100
$ws.Range("O2").Value = 0

$ws.Range("H3").Value = 3

$ws.Range("H4").Value = 4
$ws.Range("J4").Value = 0.5
$ws.Range("L4").Value = 4

$ws.Range("H5").Value = 4
$ws.Range("J5").Value = 0.5
$ws.Range("L5").Value = 4

$ws.Range("H6").Value = 30
$ws.Range("J6").Value = 10
$ws.Range("L6").Value = 50

$ws.Range("H7").Value = 40
$ws.Range("J7").Value = 5
$ws.Range("L7").Value = 40

$ws.Range("H8").Value = 40
$ws.Range("J8").Value = 5
$ws.Range("L8").Value = 40

# --- Column width for column T (20) widened and no longer auto (bestFit) ---
$ws.Columns.Item(20).ColumnWidth = 16.43

# --- Selection moved from P4 to P7 ---
$ws.Range("P7").Select()
